$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns AD/AE/AF -> Wins/Losses/Ties
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, centered, bordered) from an existing header cell (AC1)
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Re-apply the values after pasting formats (PasteSpecial formats only, so values remain,
# but set again to be safe)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-49: team record Wins=90, Losses=73, Ties=0
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 30).Value = 90   # AD
    $ws.Cells.Item($r, 31).Value = 73   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
